# Fruta / hortaliza, semanal
# Insert a new weekly record at row 67 of the "Frambuesa" (raspberry) price
# table (Vega Central Mapocho de Santiago). All existing data rows from 67
# downward shift down by one row, and the new row is populated with the
# latest weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 67..90 down to 68..91, leaving a blank row 67 for the new record.
$ws.Rows(67).Insert()

$newRow = 67

$ws.Cells.Item($newRow, 1).Value  = 9
$ws.Cells.Item($newRow, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($newRow, 3).Value  = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value  = 44636
$ws.Cells.Item($newRow, 5).Value  = 13
$ws.Cells.Item($newRow, 6).Value  = "Fruta"
$ws.Cells.Item($newRow, 7).Value  = 100101
$ws.Cells.Item($newRow, 8).Value  = "Berries"
$ws.Cells.Item($newRow, 9).Value  = 100101004
$ws.Cells.Item($newRow, 10).Value = "Frambuesa"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 450
$ws.Cells.Item($newRow, 14).Value = 8000
$ws.Cells.Item($newRow, 15).Value = 8000
$ws.Cells.Item($newRow, 16).Value = 8000
$ws.Cells.Item($newRow, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item($newRow, 18).Value = "Provincia de Linares"
$ws.Cells.Item($newRow, 19).Value = 4000
$ws.Cells.Item($newRow, 20).Value = 2
